$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain plain text so values like
# "77.30" or "1.001" are not reinterpreted as numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.868.04'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.728.42'
$ws.Range('E3').Value = '  +0.22%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9978'
$ws.Range('E4').Value = '  -0.23%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.79'
$ws.Range('E5').Value = '  -0.78%  '

$ws.Range('E6').Value = '  -0.20%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2599'
$ws.Range('E8').Value = '  -0.35%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06219'
$ws.Range('E9').Value = '  +0.44%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.734.00'
$ws.Range('E10').Value = '  +0.55%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.02'
$ws.Range('E11').Value = '  +3.39%  '

$ws.Range('E12').Value = '  -1.29%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6104'
$ws.Range('E13').Value = '  +1.86%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.495'
$ws.Range('E14').Value = '  -1.68%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.30'
$ws.Range('E15').Value = '  +0.15%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9989'
$ws.Range('E16').Value = '  -0.14%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.636.66'
$ws.Range('E17').Value = '  +0.83%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9978'
$ws.Range('E18').Value = '  -0.27%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007194'
$ws.Range('E19').Value = '  +1.12%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.45'
$ws.Range('E20').Value = '  +1.01%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.956.24'
$ws.Range('E21').Value = '  +0.56%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.437'
$ws.Range('E22').Value = '  -0.63%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.572'
$ws.Range('E23').Value = '  -0.16%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.126'
$ws.Range('E24').Value = '  -0.53%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '138.46'
$ws.Range('E25').Value = '  +0.80%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.31'
$ws.Range('E26').Value = '  +0.60%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.780'
$ws.Range('E27').Value = '  +4.80%  '

$ws.Range('E28').Value = '  -0.74%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '106.30'
$ws.Range('E29').Value = '  -0.63%  '

$ws.Range('E30').Value = '  +0.37%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07994'
$ws.Range('E31').Value = '  +0.76%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.688'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04536'
$ws.Range('E33').Value = '  +0.07%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.611'
$ws.Range('E34').Value = '  +0.29%  '

$ws.Range('E35').Value = '  +1.56%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6253'
$ws.Range('E36').Value = '  -0.02%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9315'
$ws.Range('E37').Value = '  +1.97%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.060'
$ws.Range('E38').Value = '  +5.65%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.443'
$ws.Range('E39').Value = '  +2.14%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.001'
$ws.Range('E40').Value = '  +0.01%  '

$ws.Range('E41').Value = '  +1.53%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.678'
$ws.Range('E42').Value = '  +4.45%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.70'
$ws.Range('E43').Value = '  -0.25%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3867'
$ws.Range('E44').Value = '  +0.74%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.954'
$ws.Range('E45').Value = '  +3.90%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1160'
$ws.Range('E46').Value = '  +0.36%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05390'
$ws.Range('E47').Value = '  +0.48%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.970'
$ws.Range('E48').Value = '  +3.65%  '

$ws.Range('E49').Value = '  +0.49%  '

$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.240'
$ws.Range('E50').Value = '  +0.38%  '

$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.72'
$ws.Range('E51').Value = '  +1.68%  '
